$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "index" column header (A1) to "i".
# Excel will automatically keep the shared-strings table / table1.xml column
# name in sync since the table's header cell (A1) drives the column name.
$ws.Range("A1").Value = "i"

# Re-number the "index" column values so they are 0-based instead of 1-based
# (row 2 -> 0, row 3 -> 1, ... row 503 -> 501).
$rng = $ws.Range("A2:A503")
$vals = $rng.Value2
for ($i = 1; $i -le $vals.GetLength(0); $i++) {
    $vals[$i,1] = $vals[$i,1] - 1
}
$rng.Value2 = $vals

# Narrow column A now that the values/header are shorter.
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665
